$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.209707140922546
$ws.Range("B1").Value = 2.279402256011963
$ws.Range("D1").Value = 1.430251479148865
$ws.Range("E1").Value = 0.907259464263916
